$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "3、存量资产-y" : fill in the December 2024 (column L) data
# ---------------------------------------------------------------
$wsY = $wb.Worksheets.Item("3、存量资产-y")

$wsY.Range("L2").Value  = 11585
$wsY.Range("L3").Value  = 105966
$wsY.Range("L4").Value  = 0
$wsY.Range("L5").Value  = 14792
$wsY.Range("L6").Value  = 122464
$wsY.Range("L8").Value  = 22171
$wsY.Range("L9").Value  = 5.65
$wsY.Range("L10").Value = 0
$wsY.Range("L13").Value = 20049
$wsY.Range("L14").Value = 394720
$wsY.Range("L15").Value = 32384
$wsY.Range("L16").Value = 102293
$wsY.Range("L18").Value = 135584
$wsY.Range("L19").Value = 133109
$wsY.Range("L20").Value = 22888

# L17 used to hold a shared-formula reference; rewrite it explicitly
$wsY.Range("L17").Formula = "=SUM(L2:L16)"

# Update the selected cell for this sheet
$wsY.Range("L26").Select()

# ---------------------------------------------------------------
# Sheet "3、存量资产-k" : fill in the December 2025 (column L) data
# ---------------------------------------------------------------
$wsK = $wb.Worksheets.Item("3、存量资产-k")

$wsK.Range("L2").Formula = "=473991*7.1"
$wsK.Range("L3").Formula = "=153537+12710+518+1234+307"
$wsK.Range("L4").Value   = 8176
$wsK.Range("L6").Formula = "=33429.77+58692.25+45931.2"
$wsK.Range("L7").Value   = 24142

# Update the selected cell for this sheet
$wsK.Range("L13").Select()

# ---------------------------------------------------------------
# Sheet "2、收入结余表" : progress bar + new month row
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2、收入结余表")

$ws2.Range("O12").Formula = "=(N12-5500000)/1000000"
$ws2.Range("N13").Formula = "='3、存量资产-y'!L26+'3、存量资产-k'!L13"
$ws2.Range("O13").Formula = "=(N13-5500000)/1000000"

# Update the selected cell for this sheet
$ws2.Range("P14").Select()

# ---------------------------------------------------------------
# Force a full recalculation so dependent formulas / chart caches refresh
# ---------------------------------------------------------------
$excel.CalculateFullRebuild()
